$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the short "foo/bar" requirement text with a long lorem-ipsum
# paragraph so the sheet has a wrapping scenario to test against.
$ws.Range("C2").Value = "Lorem ipsum dolor sit amet, consectetur adipiscing elit, sed do eiusmod tempor incididunt ut labore et dolore magna aliqua. Ut enim ad minim veniam, quis nostrud exercitation ullamco laboris nisi ut aliquip ex ea commodo consequat. Duis aute irure dolor in reprehenderit in voluptate velit esse cillum dolore eu fugiat nulla pariatur. Excepteur sint occaecat cupidatat non proident, sunt in culpa qui officia deserunt mollit anim id est laborum."

# Widen the "text" column so the longer content has room (36.5 -> 65).
$ws.Columns.Item(3).ColumnWidth = 64.17
